$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2203252032520325
$ws.Range("C2").Value = 0.5073170731707317
$ws.Range("J2").Value = 0.01788617886178862
$ws.Range("O2").Value = 0.0008130081300813008
$ws.Range("P2").Value = 0.1577235772357723
$ws.Range("S2").Value = 0.09593495934959349
$ws.Range("B3").Value = 0.01219512195121951
$ws.Range("C3").Value = 0.01829268292682927
$ws.Range("J3").Value = 0.02439024390243903
$ws.Range("O3").Value = 0.001524390243902439
$ws.Range("P3").Value = 0.7560975609756098
$ws.Range("S3").Value = 0.1875
$ws.Range("J4").Value = 0.07344632768361582
$ws.Range("O4").Value = 0.005649717514124294
$ws.Range("P4").Value = 0.6892655367231638
$ws.Range("S4").Value = 0.231638418079096
$ws.Range("B6").Value = 0.06793478260869565
$ws.Range("D6").Value = 0.02445652173913044
$ws.Range("F6").Value = 0.06793478260869565
$ws.Range("J6").Value = 0.2282608695652174
$ws.Range("O6").Value = 0.02309782608695652
$ws.Range("Q6").Value = 0.1861413043478261
$ws.Range("R6").Value = 0.05706521739130434
$ws.Range("S6").Value = 0.3451086956521739
$ws.Range("B7").Value = 0.1248206599713056
$ws.Range("D7").Value = 0.02439024390243903
$ws.Range("F7").Value = 0.06025824964131994
$ws.Range("J7").Value = 0.1420373027259684
$ws.Range("O7").Value = 0.01147776183644189
$ws.Range("Q7").Value = 0.1865136298421808
$ws.Range("R7").Value = 0.06743185078909612
$ws.Range("S7").Value = 0.3830703012912482
$ws.Range("B8").Value = 0.1028325123152709
$ws.Range("D8").Value = 0.02032019704433497
$ws.Range("E8").Value = 0.001231527093596059
$ws.Range("F8").Value = 0.06342364532019705
$ws.Range("J8").Value = 0.1280788177339902
$ws.Range("O8").Value = 0.01724137931034483
$ws.Range("Q8").Value = 0.1724137931034483
$ws.Range("R8").Value = 0.08312807881773399
$ws.Range("S8").Value = 0.4113300492610837
$ws.Range("B9").Value = 0.09707446808510638
$ws.Range("D9").Value = 0.01462765957446809
$ws.Range("F9").Value = 0.05319148936170213
$ws.Range("J9").Value = 0.1077127659574468
$ws.Range("O9").Value = 0.01196808510638298
$ws.Range("Q9").Value = 0.2446808510638298
$ws.Range("R9").Value = 0.07446808510638298
$ws.Range("S9").Value = 0.3962765957446808
$ws.Range("B10").Value = 0.1225693685820406
$ws.Range("D10").Value = 0.02119292112737601
$ws.Range("E10").Value = 0.0008739348918505571
$ws.Range("F10").Value = 0.05964605636880053
$ws.Range("J10").Value = 0.1238802709198165
$ws.Range("O10").Value = 0.01769718155997378
$ws.Range("Q10").Value = 0.2173913043478261
$ws.Range("R10").Value = 0.08542713567839195
$ws.Range("S10").Value = 0.351321826523924
$ws.Range("G11").Value = 0.1346153846153846
$ws.Range("J11").Value = 0.09615384615384616
$ws.Range("K11").Value = 0.1987179487179487
$ws.Range("L11").Value = 0.5512820512820513
$ws.Range("S11").Value = 0.01923076923076923
$ws.Range("G12").Value = 0.7479806138933764
$ws.Range("J12").Value = 0.1857835218093699
$ws.Range("K12").Value = 0.003231017770597738
$ws.Range("L12").Value = 0.02261712439418417
$ws.Range("S12").Value = 0.04038772213247173
$ws.Range("G13").Value = 0.6526946107784432
$ws.Range("J13").Value = 0.2754491017964072
$ws.Range("S13").Value = 0.0718562874251497
$ws.Range("G14").Value = 0.625
$ws.Range("J14").Value = 0.125
$ws.Range("S14").Value = 0.25
$ws.Range("F15").Value = 0.01783264746227709
$ws.Range("H15").Value = 0.1467764060356653
$ws.Range("I15").Value = 0.07270233196159122
$ws.Range("J15").Value = 0.3731138545953361
$ws.Range("K15").Value = 0.07407407407407407
$ws.Range("M15").Value = 0.01371742112482853
$ws.Range("O15").Value = 0.04938271604938271
$ws.Range("S15").Value = 0.252400548696845
$ws.Range("F16").Value = 0.02287166454891995
$ws.Range("H16").Value = 0.1613722998729352
$ws.Range("I16").Value = 0.06226175349428208
$ws.Range("J16").Value = 0.4129606099110547
$ws.Range("K16").Value = 0.1143583227445997
$ws.Range("M16").Value = 0.01905972045743329
$ws.Range("N16").Value = 0.006353240152477764
$ws.Range("O16").Value = 0.05844980940279543
$ws.Range("S16").Value = 0.1423125794155019
$ws.Range("F17").Value = 0.01792943898207056
$ws.Range("H17").Value = 0.183342972816657
$ws.Range("I17").Value = 0.09947946790052054
$ws.Range("J17").Value = 0.4112203585887796
$ws.Range("K17").Value = 0.09427414690572586
$ws.Range("M17").Value = 0.02486986697513013
$ws.Range("N17").Value = 0.002313475997686524
$ws.Range("O17").Value = 0.05320994794679005
$ws.Range("S17").Value = 0.1133603238866397
$ws.Range("F18").Value = 0.01337295690936107
$ws.Range("H18").Value = 0.187221396731055
$ws.Range("I18").Value = 0.09212481426448738
$ws.Range("J18").Value = 0.424962852897474
$ws.Range("K18").Value = 0.09361069836552749
$ws.Range("M18").Value = 0.02526002971768202
$ws.Range("O18").Value = 0.0549777117384844
$ws.Range("S18").Value = 0.1084695393759287
$ws.Range("F19").Value = 0.01523980277902286
$ws.Range("H19").Value = 0.2140295831465711
$ws.Range("I19").Value = 0.09188704616763783
$ws.Range("J19").Value = 0.3644105782160466
$ws.Range("K19").Value = 0.1102644554011654
$ws.Range("M19").Value = 0.01927386822052891
$ws.Range("N19").Value = 0.0006723442402510085
$ws.Range("O19").Value = 0.06432093231734648
$ws.Range("S19").Value = 0.1199013895114299
